# Update scripts with new TPM (transcripts per million) values.
#
# The "Sending cluster" x "Target cluster" communication table is refreshed
# with newly-computed statistics (columns E:T), and a second sending cluster
# ("MuSCs") is added with one row for every target cluster, extending the
# sheet from 5 data rows (A2:T6) to 10 data rows (A2:T11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (A..T) in order, matching the existing row 1 layout.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Row 2 .. 11 values for every column. Sending cluster (A) / Ligand (B) /
# Receptor (C) are fixed across all rows for this ligand-receptor pair; the
# Target cluster (D) and the numeric statistics (E..T) vary per row.
$data = @(
    @("FAPs",  "Inha", "Tgfbr3", "ECs",              3, 1,                  0.5953293333333333, 1.785988, 0.8539648188994299, 0.85396481889943,  3, 1,                  41.607769,          124.823307, 0.4674897261536314,  0.4674897261536314,  24.77032538025733,  222.932928422316,  0.3992197793321299,   0.39921977933213),
    @("FAPs",  "Inha", "Tgfbr3", "FAPs",             3, 1,                  0.5953293333333333, 1.785988, 0.8539648188994299, 0.85396481889943,  3, 1,                  40.85529333333333,  122.56588,  0.4590351838457449,  0.4590351838457449,  24.32235454327111,  218.90119088944,   0.391999897641298,    0.3919998976412981),
    @("FAPs",  "Inha", "Tgfbr3", "Inflammatory-Mac", 3, 1,                  0.5953293333333333, 1.785988, 0.8539648188994299, 0.85396481889943,  3, 1,                  0.2401453333333333, 0.720436,   0.002698185430636104, 0.002698185430636104, 0.1429655611964444,  1.286690050768,    0.002304155432630241, 0.002304155432630241),
    @("FAPs",  "Inha", "Tgfbr3", "MuSCs",            3, 1,                  0.5953293333333333, 1.785988, 0.8539648188994299, 0.85396481889943,  3, 1,                  6.288874666666666,  18.866624,  0.0706595034147231,  0.0706595034147231,  3.743951562723554,  33.69556406451199, 0.06034073003707766,  0.06034073003707767),
    @("FAPs",  "Inha", "Tgfbr3", "Resolving-Mac",    3, 1,                  0.5953293333333333, 1.785988, 0.8539648188994299, 0.85396481889943,  1, 0.3333333333333333, 0.010449,           0.031347,   0.0001174011552645203, 0.0001174011552645203, 0.006220596203999999, 0.05598536583599999, 0.0001002564562940499, 0.0001002564562940499),
    @("MuSCs", "Inha", "Tgfbr3", "ECs",              2, 0.6666666666666666, 0.1018063333333333, 0.305419, 0.1460351811005701, 0.1460351811005701, 3, 1,                  41.607769,          124.823307, 0.4674897261536314,  0.4674897261536314,  4.235934400070333,  38.123409600633,   0.06826994682150149,  0.06826994682150149),
    @("MuSCs", "Inha", "Tgfbr3", "FAPs",             2, 0.6666666666666666, 0.1018063333333333, 0.305419, 0.1460351811005701, 0.1460351811005701, 3, 1,                  40.85529333333333,  122.56588,  0.4590351838457449,  0.4590351838457449,  4.159327611524444,  37.43394850372,    0.06703528620444683,  0.06703528620444683),
    @("MuSCs", "Inha", "Tgfbr3", "Inflammatory-Mac", 2, 0.6666666666666666, 0.1018063333333333, 0.305419, 0.1460351811005701, 0.1460351811005701, 3, 1,                  0.2401453333333333, 0.720436,   0.002698185430636104, 0.002698185430636104, 0.02444831585377778, 0.220034842684,    0.0003940299980058632, 0.0003940299980058632),
    @("MuSCs", "Inha", "Tgfbr3", "MuSCs",            2, 0.6666666666666666, 0.1018063333333333, 0.305419, 0.1460351811005701, 0.1460351811005701, 3, 1,                  6.288874666666666,  18.866624,  0.0706595034147231,  0.0706595034147231,  0.6402472706062221, 5.762225435455999, 0.01031877337764544,  0.01031877337764544),
    @("MuSCs", "Inha", "Tgfbr3", "Resolving-Mac",    2, 0.6666666666666666, 0.1018063333333333, 0.305419, 0.1460351811005701, 0.1460351811005701, 1, 0.3333333333333333, 0.010449,           0.031347,   0.0001174011552645203, 0.0001174011552645203, 0.001063774377,      0.009573969393,    0.00001714469897047037, 0.00001714469897047037)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$rowNum").Value = $rowValues[$j]
    }
}
